# Update Work Week and Social Spending
# Updates GDP per Capita data for Qatar (1950-2010) and extends through 2016
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update existing rows (years 1950-2010) in column E with revised GDP per Capita values ---
$ws.Range("E2").Formula = '="48436"'
$ws.Range("E3").Formula = '="48653"'
$ws.Range("E4").Formula = '="47990"'
$ws.Range("E5").Formula = '="49284"'
$ws.Range("E6").Formula = '="51490"'
$ws.Range("E7").Formula = '="49632"'
$ws.Range("E8").Formula = '="50628"'
$ws.Range("E9").Formula = '="49973"'
$ws.Range("E10").Formula = '="50042"'
$ws.Range("E11").Formula = '="52436"'
$ws.Range("E12").Formula = '="52299"'
$ws.Range("E13").Formula = '="48516"'
$ws.Range("E14").Formula = '="46304"'
$ws.Range("E15").Formula = '="44912"'
$ws.Range("E16").Formula = '="42186"'
$ws.Range("E17").Formula = '="41101"'
$ws.Range("E18").Formula = '="50660"'
$ws.Range("E19").Formula = '="55566"'
$ws.Range("E20").Formula = '="58010"'
$ws.Range("E21").Formula = '="56236"'
$ws.Range("E22").Formula = '="51921"'
$ws.Range("E23").Formula = '="59731"'
$ws.Range("E24").Formula = '="62416"'
$ws.Range("E25").Formula = '="68407"'
$ws.Range("E26").Formula = '="57594"'
$ws.Range("E27").Formula = '="54866"'
$ws.Range("E28").Formula = '="55169"'
$ws.Range("E29").Formula = '="45999"'
$ws.Range("E30").Formula = '="47072"'
$ws.Range("E31").Formula = '="45806"'
$ws.Range("E32").Formula = '="45860"'
$ws.Range("E33").Formula = '="37307"'
$ws.Range("E34").Formula = '="29046"'
$ws.Range("E35").Formula = '="23146"'
$ws.Range("E36").Formula = '="20354"'
$ws.Range("E37").Formula = '="16558"'
$ws.Range("E38").Formula = '="12881"'
$ws.Range("E39").Formula = '="12819"'
$ws.Range("E40").Formula = '="12492"'
$ws.Range("E41").Formula = '="12137"'
$ws.Range("E42").Formula = '="11705"'
$ws.Range("E43").Formula = '="12417.6427075336"'
$ws.Range("E44").Formula = '="15042.1231591104"'
$ws.Range("E45").Formula = '="16221.0455709903"'
$ws.Range("E46").Formula = '="17955.3148177473"'
$ws.Range("E47").Formula = '="19959.7030599908"'
$ws.Range("E48").Formula = '="22395.2638874842"'
$ws.Range("E49").Formula = '="30989.1497275617"'
$ws.Range("E50").Formula = '="36448.0565388444"'
$ws.Range("E51").Formula = '="40135.174729466"'
$ws.Range("E52").Formula = '="45788.122704193"'
$ws.Range("E53").Formula = '="50525.7058393558"'
$ws.Range("E54").Formula = '="57541.103389298"'
$ws.Range("E55").Formula = '="62234.2138242042"'
$ws.Range("E56").Formula = '="74372.9325180593"'
$ws.Range("E57").Formula = '="76808.4200141651"'
$ws.Range("E58").Formula = '="90120.917983189"'
$ws.Range("E59").Formula = '="97912.6904920719"'
$ws.Range("E60").Formula = '="107402.056105853"'
$ws.Range("E61").Formula = '="115283.416602271"'
$ws.Range("E62").Formula = '="134802.78185179"'

# Convert the formulas to static text values (preserves text type without adding cell styles)
$ws.Range("E2:E62").Copy()
$ws.Range("E2:E62").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Append new rows for years 2011-2016 ---
$ws.Range("A63").Value = 634
$ws.Range("B63").Formula = '="Qatar"'
$ws.Range("C63").Formula = '="GDP per Capita"'
$ws.Range("D63").Value = 2011
$ws.Range("E63").Formula = '="155533"'
$ws.Range("A64").Value = 634
$ws.Range("B64").Formula = '="Qatar"'
$ws.Range("C64").Formula = '="GDP per Capita"'
$ws.Range("D64").Value = 2012
$ws.Range("E64").Formula = '="153922"'
$ws.Range("A65").Value = 634
$ws.Range("B65").Formula = '="Qatar"'
$ws.Range("C65").Formula = '="GDP per Capita"'
$ws.Range("D65").Value = 2013
$ws.Range("E65").Formula = '="154159"'
$ws.Range("A66").Value = 634
$ws.Range("B66").Formula = '="Qatar"'
$ws.Range("C66").Formula = '="GDP per Capita"'
$ws.Range("D66").Value = 2014
$ws.Range("E66").Formula = '="155069"'
$ws.Range("A67").Value = 634
$ws.Range("B67").Formula = '="Qatar"'
$ws.Range("C67").Formula = '="GDP per Capita"'
$ws.Range("D67").Value = 2015
$ws.Range("E67").Formula = '="156029"'
$ws.Range("A68").Value = 634
$ws.Range("B68").Formula = '="Qatar"'
$ws.Range("C68").Formula = '="GDP per Capita"'
$ws.Range("D68").Value = 2016
$ws.Range("E68").Formula = '="156299"'

# Convert new-row formulas to static text values
$ws.Range("B63:C68").Copy()
$ws.Range("B63:C68").PasteSpecial(-4163)
$ws.Range("E63:E68").Copy()
$ws.Range("E63:E68").PasteSpecial(-4163)
$excel.CutCopyMode = 0
